# Auto-generated script applying the "Knärot" section addition and date update
# described by the commit diff.
$d = $word.ActiveDocument

# --- Part 1: insert the new "Knärot" content block after the last paragraph
#     ("BILAGA 1 - Fridlysta arter"), before the final section break.
#
# Phase 1a: create every paragraph with its plain (un-italicised) text and the
# correct paragraph style, while remembering where each paragraph begins/ends so
# that italics can be applied afterwards. (Applying italics immediately, before
# later paragraphs are inserted, causes the emphasis to bleed into the text that
# is typed afterwards, so all insertion happens first and all italicising last.)
$lastIndex = $d.Paragraphs.Count
$insertAfterRange = $d.Paragraphs.Item($lastIndex).Range

$italicJobs = @()

# Paragraph 1: style=Heading 1
$insertAfterRange.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Heading 1"
$pStart = $p.Range.Start
$p.Range.InsertAfter("Knärot – ekologi samt krav på livsmiljön")
$pEnd = $p.Range.End
$insertAfterRange = $p.Range

# Paragraph 2: style=Normal
$insertAfterRange.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"
$pStart = $p.Range.Start
$p.Range.InsertAfter("Knärot är fridlyst enligt 8 och 15 §§ artskyddsförordningen och klassad som sårbar (VU) enligt rödlistan 2020. Knärot är beroende av hög och jämn luftfuktighet i gamla, ostörda skogsmiljöer och är känslig för snabba förändringar av ljus-/vindförhållanden eller uttorkning. På grund av ett alltför intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 åren och i framtiden bedöms minskningstakten uppgå till 30 (20-40) %. Till följd av att arten har en dokumenterat högre minskningstakt iförhållande till sin generationstid än vad som tidigare varit känt (data från Riksskogstaxeringen) höjdes den till hotkategori sårbar (VU) i rödlistan 2020 (Artdatabanken, 2021).")
$pEnd = $p.Range.End
$insertAfterRange = $p.Range

# Paragraph 3: style=Normal
$insertAfterRange.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"
$pStart = $p.Range.Start
$p.Range.InsertAfter("Samuel Johnsons doktorsavhandling “Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“ (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: “Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” Vidare “More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”")
$pEnd = $p.Range.End
$italicJobs += , @($pStart, $pEnd, "“Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“")
$italicJobs += , @($pStart, $pEnd, "“Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” ")
$italicJobs += , @($pStart, $pEnd, "“More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”")
$insertAfterRange = $p.Range

# Paragraph 4: style=Normal
$insertAfterRange.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"
$pStart = $p.Range.Start
$p.Range.InsertAfter("Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: “In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”")
$pEnd = $p.Range.End
$italicJobs += , @($pStart, $pEnd, "“In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”")
$insertAfterRange = $p.Range

# Paragraph 5: style=Normal
$insertAfterRange.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"
$pStart = $p.Range.Start
$p.Range.InsertAfter("En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkidén knärots skyddsbehov. I uppsatsen berörs problemet med uttorkning för växter, bl.a. för knärot, ett problem som blivit accentuerat på grund av den pågående klimatförändringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen undersöks områden med tre olika avstånd från kalhyggeskant med avseende på skydd bl.a. för knärot. Det första området har avstånd upp till 20 m från hyggeskant (Strong edge effect), det andra 20 – 40 m från hyggeskant (Weak edge effect) och det tredje avser större avstånd från hyggeskant, där kanteffekten anses vara försumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt på känsliga och rödlistade skogsarter vid de kortare avstånden till hyggeskant, medan effekt av uttorkning inte konstaterades på större avstånd (Interior). För orkidén knärot fann man en rik förekomst (upp till 0,06 dm2/m2) på stort avstånd från hyggeskant (Interior), medan förekomsten var liten eller närmast försumbar i de områden som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet påpekar att de allt oftare förekommande torra somrarna ger ytterligare skäl att utöka skyddsavståndet från hyggen till den fuktkrävande arten knärot (Koelmeijer m.fl., 2022).")
$pEnd = $p.Range.End
$insertAfterRange = $p.Range

# Paragraph 6: style=Normal
$insertAfterRange.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"
$pStart = $p.Range.Start
$p.Range.InsertAfter("Även Skogsstyrelsens egen vägledning för hänsyn till knärot ligger i linje med ovanstående forskningsstudier. Av vägledningen framgår det att för med hög sannolikhet kunna bevara befintliga förekomster krävs relativt stora avsättningar av uppvuxen skog med slutet och relativt tätt kronskikt. Som riktlinje kan krävas ett avstånd på 50 meter in från brynet för att vidmakthålla ett fungerande mikroklimat. Detta innebär att fristående hänsynsytor för många arter (kärlväxter, lavar och mossor) kan behöva ha en area överstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) för att bibehålla lokalklimatet. Även ganska små förändringar i form av förändrade ljus- och fuktighetsförhållanden, till exempel till följd av gallring, kan leda till att arten försvinner till följd av konkurrens med mera ljuskrävande och snabbväxande arter (Skogsstyrelsen, 2022).")
$pEnd = $p.Range.End
$insertAfterRange = $p.Range

# Paragraph 7: style=Heading 2
$insertAfterRange.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Heading 2"
$pStart = $p.Range.Start
$p.Range.InsertAfter("Referenser - knärot")
$pEnd = $p.Range.End
$insertAfterRange = $p.Range

# Paragraph 8: style=Normal
$insertAfterRange.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"
$pStart = $p.Range.Start
$p.Range.InsertAfter("de Graaf M & Roberts M.R., 2009. Short-term response of the herbaceous layer within leave patches after harvest. Forest Ecology and Management 257, 1014-1025")
$pEnd = $p.Range.End
$italicJobs += , @($pStart, $pEnd, "Short-term response of the herbaceous layer within leave patches after harvest. ")
$insertAfterRange = $p.Range

# Paragraph 9: style=Normal
$insertAfterRange.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"
$pStart = $p.Range.Start
$p.Range.InsertAfter("Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. Ecological Applications, 22, 2049-2064 ")
$pEnd = $p.Range.End
$italicJobs += , @($pStart, $pEnd, "Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. ")
$insertAfterRange = $p.Range

# Paragraph 10: style=Normal
$insertAfterRange.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"
$pStart = $p.Range.Start
$p.Range.InsertAfter("Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. Interactive effects of drought and edge exposure on old-growth forest understory species. Landscape Ecology, 37, sid 1839-1853")
$pEnd = $p.Range.End
$italicJobs += , @($pStart, $pEnd, "Interactive effects of drought and edge exposure on old-growth forest understory species. ")
$insertAfterRange = $p.Range

# Paragraph 11: style=Normal
$insertAfterRange.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"
$pStart = $p.Range.Start
$p.Range.InsertAfter("Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. Biological legacies buffer local species extinction after logging. Journal of Applied Ecology. 51, 53-62.")
$pEnd = $p.Range.End
$italicJobs += , @($pStart, $pEnd, "Biological legacies buffer local species extinction after logging. ")
$insertAfterRange = $p.Range

# Paragraph 12: style=Normal
$insertAfterRange.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"
$pStart = $p.Range.Start
$p.Range.InsertAfter("Skogsstyrelsen, 2022. Vägledning för hänsyn till knärot. https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/")
$pEnd = $p.Range.End
$italicJobs += , @($pStart, $pEnd, "Vägledning för hänsyn till knärot. ")
$insertAfterRange = $p.Range

# Paragraph 13: style=Normal
$insertAfterRange.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"
$pStart = $p.Range.Start
$p.Range.InsertAfter("SLU Artdatabanken, 2021. Artfaktablad. Naturvård – artfakta. SLU Artdatabanken, Uppsala ")
$pEnd = $p.Range.End
$italicJobs += , @($pStart, $pEnd, "Artfaktablad. Naturvård – artfakta. ")
$insertAfterRange = $p.Range

# Phase 1b: now that all paragraphs/text exist, apply italics to the recorded
# sub-ranges. Doing this last avoids the formatting leaking into later inserts.
foreach ($job in $italicJobs) {
    $searchRange = $d.Range($job[0], $job[1])
    $found = $searchRange.Find.Execute($job[2], $false, $false, $false, $false, $false,
                                        $true, 1, $false, "", 0)
    if (-not $found) { throw "Could not find italic run text: $($job[2])" }
    $searchRange.Font.Italic = 1
}

# --- Part 2: update the header date from 2023-09-13 to 2023-09-15
foreach ($sec in $d.Sections) {
    foreach ($hf in $sec.Headers) {
        if ($hf.Exists) {
            $hf.Range.Find.Execute("2023-09-13", $true, $false, $false, $false, $false,
                                    $true, 1, $false, "2023-09-15", 2) | Out-Null
        }
    }
    foreach ($ft in $sec.Footers) {
        if ($ft.Exists) {
            $ft.Range.Find.Execute("2023-09-13", $true, $false, $false, $false, $false,
                                    $true, 1, $false, "2023-09-15", 2) | Out-Null
        }
    }
}
$d.Content.Find.Execute("2023-09-13", $true, $false, $false, $false, $false,
                        $true, 1, $false, "2023-09-15", 2) | Out-Null

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count